$wb = $excel.ActiveWorkbook

# --- "About" sheet (sheet 1) ---
$about = $wb.Worksheets.Item("About")

# Values are entered in this specific order to reproduce the author's
# original shared-string table ordering.
$about.Range("A14").Value = "input data to avoid runtime crashes:"

$about.Range("A15").Value = "elec/CCAMC"
$about.Range("A15").Font.Bold = $true

$about.Range("C15").Value = "capacity costs for wind and solar in start year"

$about.Range("A13").Value = "If advancing the initial time in this variable, you must update the following variable's"

$about.Range("A18").Value = "plcy-schd/FoPITY"
$about.Range("A18").Font.Bold = $true

$about.Range("C18").Value = "polcy implementation schedule"

$about.Range("A17").Value = "Other variables whose input data it would be wise to update:"

$about.Range("A19").Value = "trans/SYVbT"
$about.Range("A19").Font.Bold = $true

$about.Range("A20").Value = "other ""Start Year"" variables"
$about.Range("A20").Font.Bold = $true

# --- "IT" sheet (sheet 2) ---
$it = $wb.Worksheets.Item("IT")
$it.Range("B2").Value = 2020

# Selected cell / active sheet
$it.Range("B3").Select()
$it.Activate()
